# Auto-generated Excel COM-interop edit script
# Applies numeric updates to the Sheets workbook per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2549.125
$ws.Range("I100").Value = 1580
$ws.Range("J100").Value = 3905.9
$ws.Range("K100").Value = 1580
$ws.Range("L100").Value = 3905.9
$ws.Range("M100").Value = -1039
$ws.Range("N100").Value = -4987.9

$ws.Range("H129").Value = 924.40424
$ws.Range("I129").Value = 253.8
$ws.Range("J129").Value = 1105.6487
$ws.Range("K129").Value = 761.4000000000001
$ws.Range("L129").Value = 3316.9461
$ws.Range("M129").Value = 4238.6
$ws.Range("N129").Value = -13316.9461

$ws.Range("H132").Value = 1042.9181
$ws.Range("I132").Value = 998.4386
$ws.Range("J132").Value = 1676.75
$ws.Range("K132").Value = 2995.3158
$ws.Range("L132").Value = 5030.25
$ws.Range("M132").Value = -465.3157999999999
$ws.Range("N132").Value = -10090.25

$ws.Range("H133").Value = 49141.8
$ws.Range("J133").Value = 48750
$ws.Range("L133").Value = 48750
$ws.Range("N133").Value = -58870

$ws.Range("H138").Value = 2885.6177
$ws.Range("I138").Value = 1977.5
$ws.Range("J138").Value = 4035.9
$ws.Range("K138").Value = 5932.5
$ws.Range("L138").Value = 12107.7
$ws.Range("M138").Value = -792.5
$ws.Range("N138").Value = -22387.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 144777.34
$ws.Range("I32").Value = 5025.4653
$ws.Range("J32").Value = 1237382.9
$ws.Range("K32").Value = 5025.4653
$ws.Range("L32").Value = 1237382.9
$ws.Range("M32").Value = -4738.4653
$ws.Range("N32").Value = -1237956.9

$ws.Range("H37").Value = 142862860
$ws.Range("I37").Value = 250004000
$ws.Range("J37").Value = 8001
$ws.Range("K37").Value = 250004000
$ws.Range("L37").Value = 8001
$ws.Range("M37").Value = -250003727
$ws.Range("N37").Value = -8547

$ws.Range("H132").Value = 2120
$ws.Range("I132").Value = 1769.6097
$ws.Range("J132").Value = 3915.75
$ws.Range("K132").Value = 5308.8291
$ws.Range("L132").Value = 11747.25
$ws.Range("M132").Value = -2778.8291
$ws.Range("N132").Value = -16807.25

$ws.Range("H133").Value = 75026.10000000001
$ws.Range("J133").Value = 75026.10000000001
$ws.Range("L133").Value = 75026.10000000001
$ws.Range("N133").Value = -80086.10000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2936.8276
$ws.Range("I20").Value = 2403.25
$ws.Range("J20").Value = 3593.5386
$ws.Range("K20").Value = 2403.25
$ws.Range("L20").Value = 3593.5386
$ws.Range("M20").Value = -2156.25
$ws.Range("N20").Value = -4087.5386

$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -327
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H22").Value = 413.5
$ws.Range("I22").Value = 225.25
$ws.Range("J22").Value = 790
$ws.Range("K22").Value = 225.25
$ws.Range("L22").Value = 790
$ws.Range("M22").Value = 124.75
$ws.Range("N22").Value = -1490

$ws.Range("H31").Value = 9833.791999999999
$ws.Range("I31").Value = 3218.743
$ws.Range("J31").Value = 15346.333
$ws.Range("K31").Value = 3218.743
$ws.Range("L31").Value = 15346.333
$ws.Range("M31").Value = -2923.743
$ws.Range("N31").Value = -15936.333

$ws.Range("H34").Value = 9833.791999999999
$ws.Range("I34").Value = 3218.743
$ws.Range("J34").Value = 15346.333
$ws.Range("K34").Value = 3218.743
$ws.Range("L34").Value = 15346.333
$ws.Range("M34").Value = -3016.743
$ws.Range("N34").Value = -15750.333

$ws.Range("H60").Value = 7899.2
$ws.Range("I60").Value = 7093
$ws.Range("J60").Value = 8100.75
$ws.Range("K60").Value = 7093
$ws.Range("L60").Value = 8100.75
$ws.Range("M60").Value = -6582
$ws.Range("N60").Value = -9122.75

$ws.Range("H68").Value = 17073.75
$ws.Range("J68").Value = 17073.75
$ws.Range("L68").Value = 17073.75
$ws.Range("N68").Value = -18571.75

$ws.Range("H71").Value = 17073.75
$ws.Range("J71").Value = 17073.75
$ws.Range("L71").Value = 51221.25
$ws.Range("N71").Value = -58709.25

$ws.Range("H74").Value = 18778
$ws.Range("J74").Value = 18778
$ws.Range("L74").Value = 18778
$ws.Range("N74").Value = -20526

$ws.Range("H77").Value = 18778
$ws.Range("J77").Value = 18778
$ws.Range("L77").Value = 56334
$ws.Range("N77").Value = -65070

$ws.Range("H132").Value = 5018.5
$ws.Range("I132").Value = 3999.5
$ws.Range("J132").Value = 5528
$ws.Range("K132").Value = 11998.5
$ws.Range("L132").Value = 16584
$ws.Range("M132").Value = -9468.5
$ws.Range("N132").Value = -21644

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1967
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 1960.4
$ws.Range("K25").Value = 6000
$ws.Range("L25").Value = 5881.200000000001
$ws.Range("M25").Value = -5831
$ws.Range("N25").Value = -6219.200000000001

$ws.Range("H30").Value = 1967
$ws.Range("I30").Value = 2000
$ws.Range("J30").Value = 1960.4
$ws.Range("K30").Value = 6000
$ws.Range("L30").Value = 5881.200000000001
$ws.Range("M30").Value = -5898
$ws.Range("N30").Value = -6085.200000000001

$ws.Range("H131").Value = 11765524
$ws.Range("I131").Value = 20000796
$ws.Range("J131").Value = 8334160.5
$ws.Range("K131").Value = 60002388
$ws.Range("L131").Value = 25002481.5
$ws.Range("M131").Value = -59997348
$ws.Range("N131").Value = -25012561.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2234.976
$ws.Range("I122").Value = 1894.3
$ws.Range("K122").Value = 5682.9
$ws.Range("M122").Value = -3232.9

$ws.Range("H132").Value = 2682.2285
$ws.Range("I132").Value = 2132.8845
$ws.Range("J132").Value = 4269.222
$ws.Range("K132").Value = 6398.6535
$ws.Range("L132").Value = 12807.666
$ws.Range("M132").Value = -3868.6535
$ws.Range("N132").Value = -17867.666

$ws.Range("H133").Value = 38399.25
$ws.Range("J133").Value = 38399.25
$ws.Range("L133").Value = 38399.25
$ws.Range("N133").Value = -48519.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H132").Value = 5479.524
$ws.Range("I132").Value = 6264.5713
$ws.Range("J132").Value = 3909.4285
$ws.Range("K132").Value = 18793.7139
$ws.Range("L132").Value = 11728.2855
$ws.Range("M132").Value = -16263.7139
$ws.Range("N132").Value = -16788.2855

$ws.Range("H133").Value = 80103.06
$ws.Range("J133").Value = 80103.06
$ws.Range("L133").Value = 80103.06
$ws.Range("N133").Value = -85163.06

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 100013210
$ws.Range("I81").Value = 2869
$ws.Range("J81").Value = 333370660
$ws.Range("K81").Value = 5738
$ws.Range("L81").Value = 666741320
$ws.Range("M81").Value = -4677
$ws.Range("N81").Value = -666743442

$ws.Range("H84").Value = 100013210
$ws.Range("I84").Value = 2869
$ws.Range("J84").Value = 333370660
$ws.Range("K84").Value = 28690
$ws.Range("L84").Value = 3333706600
$ws.Range("M84").Value = -23386
$ws.Range("N84").Value = -3333717208

$ws.Range("H132").Value = 39478676
$ws.Range("I132").Value = 71430730
$ws.Range("J132").Value = 8499.294
$ws.Range("K132").Value = 214292190
$ws.Range("L132").Value = 25497.882
$ws.Range("M132").Value = -214289660
$ws.Range("N132").Value = -30557.882

$ws.Range("H133").Value = 28000
$ws.Range("J133").Value = 28000
$ws.Range("L133").Value = 28000
$ws.Range("N133").Value = -38120
